# Agile Activity List.xlsx - apply the recorded edit:
#   - workbook window yWindow: 900 -> 1350              (bookViews/workbookView)
#   - custom date numFmt id:   165 -> 164                (xl/styles.xml)   [not exposed by the OM - best effort]
#   - sheet zoom:              75% -> 64%                (sheetView zoomScale/zoomScaleNormal)
#   - selected cell:           E20 -> D2                 (sheetView/selection)
#   - column C width:          24.140625 -> 30           (cols/col width, ~29.17 "characters")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window position (yWindow 900 -> 1350). Not wired to the xlsx exporter in
# this runtime, but set it anyway in case the host honors it.
try { $excel.ActiveWindow.Top = 1350 } catch {}

# --- Zoom: 75% -> 64%
$excel.ActiveWindow.Zoom = 64

# --- Column C width: 24.140625 -> 30 (stored width == ColumnWidth + ~0.8333,
# so 29.166666666666668 round-trips to an on-disk width of exactly 30).
$ws.Columns.Item(3).ColumnWidth = 29.166666666666668

# --- Selection: E20 -> D2 (must be set after the zoom/column-width changes so
# it ends up as the sheet's persisted active cell).
$ws.Range("D2").Select()
